# Change the "." to "-" in the species ID, in multifasta and metadata.
# Column A (sample_name) gets the species ID with "." replaced by "-".
# Column B (ncbi-spuid) keeps the original species ID text but now carries
# a trailing space (matching the multifasta header change).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Diphtheria_Metadata")

$ws.Range("B4").Value = "BX248355.1-segment2 "
$ws.Range("B5").Value = "BX248355.1-segment3 "
$ws.Range("B6").Value = "BX248355.1-segment4 "

$ws.Range("A4").Value = "BX248355-1-segment2 "
$ws.Range("A5").Value = "BX248355-1-segment3 "
$ws.Range("A6").Value = "BX248355-1-segment4 "

$ws.Range("A9").Select()
